$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.659.89'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.829.98'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '410.92'
$ws.Range('E5').Value = '  -1.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.09'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.820.27'
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.613'
$ws.Range('E8').Value = '  -5.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.731'
$ws.Range('E10').Value = '  -5.68%  '
$ws.Range('E11').Value = '  -8.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000365'
$ws.Range('E12').Value = '  -9.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.88'
$ws.Range('E13').Value = '  -5.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.417.11'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.04'
$ws.Range('E15').Value = '  -5.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.38'
$ws.Range('E16').Value = '  +14.32%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.833.19'
$ws.Range('E17').Value = '  +3.94%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.138'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  -5.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.091.85'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('E21').Value = '  -6.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '410.89'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.51'
$ws.Range('E23').Value = '  -10.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.20'
$ws.Range('E24').Value = '  -5.06%  '
$ws.Range('E25').Value = '  -3.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '36.68'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.70'
$ws.Range('E27').Value = '  +11.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.15'
$ws.Range('E28').Value = '  -5.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.52'
$ws.Range('E29').Value = '  -6.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '696.77'
$ws.Range('E30').Value = '  +6.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.44'
$ws.Range('E31').Value = '  -2.22%  '
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('E33').Value = '  +0.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.14'
$ws.Range('E34').Value = '  -2.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.151'
$ws.Range('E35').Value = '  -7.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.54'
$ws.Range('E36').Value = '  -8.23%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0798'
$ws.Range('E38').Value = '  +7.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '54.76'
$ws.Range('E39').Value = '  -4.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.17'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0457'
$ws.Range('E41').Value = '  -7.74%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.136'
$ws.Range('E43').Value = '  -8.88%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '148.93'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.50'
$ws.Range('E45').Value = '  +2.76%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.15'
$ws.Range('E46').Value = '  -8.16%  '
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.26'
$ws.Range('E49').Value = '  -10.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.79'
$ws.Range('E50').Value = '  -3.85%  '
$ws.Range('E51').Value = '  -4.42%  '
